# Update the "Generate Date" / handoff / handback timestamps in the
# handback-status report, reflecting a newer report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-01 00:49:36"

# --- zh-cn sheet ---
# H4 = Correspond Handoff Datetime, K4 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-01 00:49:31"
$wsZhCn.Range("K4").Value = "2016-09-01 00:50:12"

# --- de-de sheet ---
# H4 = Correspond Handoff Datetime (same value/string as Overview G4)
# K4 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-01 00:49:36"
$wsDeDe.Range("K4").Value = "2016-09-01 00:50:20"
